$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''38.743.24'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Value = '''2.101.16'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E4').Value = '''  +0.00%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''227.88'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  -0.44%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('E6').Value = '''  +0.43%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = '''62.28'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '''  +1.54%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = '''  -0.02%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = '''  +2.23%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''0.0840'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  -0.30%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = '''  -1.03%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''15.76'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  +6.43%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = '''2.412.68'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  +0.57%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('E14').Value = '''  -1.27%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''0.809'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  +3.45%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('E16').Value = '''  +1.24%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = '''2.106.87'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  +0.71%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''38.735.14'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  +0.54%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = '''71.80'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  +1.06%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = '''  +0.77%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = '''  +0.60%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''227.64'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D24').Value = '''2.35'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  -3.33%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''2.31'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  -0.49%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = '''9.64'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  +2.02%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = '''172.30'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''  +0.70%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = '''  +3.61%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = '''  +3.32%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''19.35'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  +1.14%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = '''  +9.05%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = '''  +0.45%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = '''  +1.19%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = '''4.75'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  -0.66%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = '''7.01'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '''  +7.40%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''0.0619'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  +1.91%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = '''2.39'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  +0.36%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = '''3.57'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  -0.24%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = '''0.999'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  +0.01%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = '''  -2.59%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''102.88'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  +2.96%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = '''  +4.02%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = '''1.529.47'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  -1.03%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = '''1.19'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  +5.98%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = '''  -0.93%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = '''7.75'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  +0.69%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = '''0.0912'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  -0.15%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = '''4.14'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  -0.55%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = '''  +1.71%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = '''  -0.84%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''2.298.83'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  +0.47%  '
$ws.Range('E51').Style = 'Normal'
